# Switch to author full name from initials with all downstream fixes.
#
# The sheets "M", "A" and "Q" each have an "author" column (column B) that
# currently holds the initials "mz" on every data row. Replace it with the
# full author name, "Maja Zaloznik", everywhere, then leave behind the
# selection state the author would have ended up with after editing each
# sheet in turn (M, then A, then Q last).

$wb = $excel.ActiveWorkbook

$oldInitials = "mz"
$newName = "Maja Založnik"

function Replace-AuthorColumn($ws) {
    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 2)
        if ($cell.Text -eq $oldInitials) {
            $cell.Value = $newName
        }
    }
}

# --- Sheet "M" ---
$wsM = $wb.Worksheets.Item("M")
Replace-AuthorColumn $wsM
$wsM.Activate()
$wsM.Range("B3:B8").Select()

# --- Sheet "A" ---
$wsA = $wb.Worksheets.Item("A")
Replace-AuthorColumn $wsA
$wsA.Activate()
$wsA.Range("B2:B6").Select()

# --- Sheet "Q" (ends up the active sheet / last selection) ---
$wsQ = $wb.Worksheets.Item("Q")
Replace-AuthorColumn $wsQ
$wsQ.Activate()
$wsQ.Range("B2:B5").Select()
